$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (scheduled GitHub Actions scrape).
# Column D ("Price") and column E ("Volume(1h)") are plain-text cells in the
# source data. Any price that looks like a plain number is pinned to Text
# format ("@") before being written so Excel keeps writing it verbatim instead
# of reinterpreting it as a number (e.g. "8.00" -> 8, "158.20" -> 158.2,
# "0.0000190" -> 1.9E-05). Prices that already contain multiple "." separators
# (e.g. "2.644.35") can never be parsed as a number, so no pinning is needed there.

$ws.Range("D2").Value = '68.207.67'
$ws.Range("E2").Value = '  -0.85%  '

$ws.Range("D3").Value = '2.644.35'
$ws.Range("E3").Value = '  -0.33%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.15'
$ws.Range("E5").Value = '  -0.75%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.14'
$ws.Range("E6").Value = '  +0.73%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  -0.44%  '

$ws.Range("E9").Value = '  +5.53%  '

$ws.Range("E10").Value = '  -0.79%  '

$ws.Range("E11").Value = '  +0.17%  '

$ws.Range("E12").Value = '  +0.81%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.11'
$ws.Range("E13").Value = '  +1.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000190'
$ws.Range("E14").Value = '  +0.81%  '

$ws.Range("E15").Value = '  -0.15%  '

$ws.Range("D16").Value = '68.144.64'
$ws.Range("E16").Value = '  -0.76%  '

$ws.Range("D17").Value = '2.635.38'
$ws.Range("E17").Value = '  -1.08%  '

$ws.Range("E18").Value = '  -0.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '362.83'
$ws.Range("E19").Value = '  -1.21%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.43'
$ws.Range("E20").Value = '  -0.82%  '

$ws.Range("E21").Value = '  +3.16%  '

$ws.Range("E22").Value = '  -0.84%  '

$ws.Range("E23").Value = '  -1.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '75.04'
$ws.Range("E24").Value = '  +2.46%  '

$ws.Range("E25").Value = '  +0.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.77'
$ws.Range("E26").Value = '  -1.19%  '

$ws.Range("E28").Value = '  +0.39%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.16%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '557.03'
$ws.Range("E30").Value = '  -3.77%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.00'
$ws.Range("E31").Value = '  +0.22%  '

$ws.Range("E32").Value = '  -0.18%  '

$ws.Range("E33").Value = '  +0.06%  '

$ws.Range("E34").Value = '  +0.88%  '

$ws.Range("E35").Value = '  +0.05%  '

$ws.Range("E36").Value = '  +1.88%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '160.95'
$ws.Range("E37").Value = '  +0.14%  '

$ws.Range("E38").Value = '  +0.79%  '

$ws.Range("E39").Value = '  +1.29%  '

$ws.Range("E40").Value = '  -1.90%  '

$ws.Range("E41").Value = '  -0.36%  '

$ws.Range("E42").Value = '  +5.22%  '

$ws.Range("E43").Value = '  +0.78%  '

$ws.Range("E44").Value = '  -0.82%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.31'
$ws.Range("E46").Value = '  -0.74%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '158.20'
$ws.Range("E47").Value = '  +1.18%  '

$ws.Range("E48").Value = '  -0.02%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '21.97'
$ws.Range("E49").Value = '  -0.10%  '

$ws.Range("E50").Value = '  -0.67%  '

$ws.Range("E51").Value = '  +0.77%  '
